$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Lasso Regression - Untuned" section (merged F1:G1)
$ws.Range("F1:G1").Merge()
$ws.Range("F1").Value = "Lasso Regression - Untuned"

# Sub-headers matching the OLS section
$ws.Range("F2").Value = "(1) Market Value (euros)"
$ws.Range("G2").Value = "(2) Natural Log of Market Value"

# Data values for the Lasso Regression - Untuned columns
$ws.Range("F3").Value = 0.3584
$ws.Range("G3").Value = 0.51549999999999996

$ws.Range("F4").Value = 0.41189999999999999
$ws.Range("G4").Value = 0.42580000000000001

$ws.Range("F5").Value = 0.64180000000000004
$ws.Range("G5").Value = 0.65259999999999996

$ws.Range("F6").Value = 1.4084000000000001
$ws.Range("G6").Value = 0.033300000000000003

$ws.Range("F7").Value = 0.25829999999999997
$ws.Range("G7").Value = 0.4642

$ws.Range("F8").Value = 0.64410000000000001
$ws.Range("G8").Value = 0.72419999999999995

# Match formatting of the corresponding OLS cells for the new columns
$ws.Range("D1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("F3:F7").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("G3:G7").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("G8").PasteSpecial(-4122)

# Autofit the new columns like the existing bestFit columns
$ws.Columns("F:G").AutoFit()

# Update selection/view to reflect where the user left off
$ws.Range("F1:G1").Select()
